# Updated hw 3 and schedule
#
# The schedule table on Sheet1 is shifted: the "TH" class meeting that used
# to be listed under row 7 (with its HW-03 related links in G7:I7) has its
# day label moved into row 5 intact, row 5's own day label (C5) is updated
# from "W" to "TH", and the G/H/I link cells that belonged to row 7 actually
# belong with row 6's lesson (05-slr-conditions / ae-05-conditions / hw-03),
# so they are relocated from G7:I7 up into G6:I6, leaving G7:I7 empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1) Row 5's day-of-week label changes from "W" to "TH".
$ws.Range("C5").Value = "TH"

# 2) Relocate the slides/ae/hw-03 links that were mistakenly on row 7
#    (G7:I7) up onto row 6 (G6:I6), preserving their formatting, and
#    remove them from row 7 entirely (including formatting), matching the
#    row 6 / row 7 cell layout in the updated schedule.
$ws.Range("G7:I7").Cut($ws.Range("G6:I6"))
$ws.Range("G7:I7").Clear()

# 3) Update the saved view/selection for the sheet to match where the
#    author left the cursor after editing.
$ws.Range("E19").Select()
